# "update claim response doccument definition"
#
# The "Elements" sheet lists FHIR Composition.section.entry rows. Previously
# there was a base "Composition.section.entry" row (row 53) followed by two
# sliced child rows: "claim" (row 54) and "signature" (row 55). The edit
# collapses this into a single "Composition.section.entry" row that directly
# carries the claim-reference details (what used to live in the "claim"
# slice row) and drops the "signature" slice row entirely. The following
# "Composition.section.emptyReason" and "Composition.section.section" rows
# shift up by two rows as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Elements")

# Drop the generic base "entry" row (row 53) -- the "claim" slice row that
# follows shifts up into its place, bringing the claim-specific Type/Short/
# Definition/Condition(s)/Mapping values with it.
$ws.Rows(53).Delete()

# The "signature" slice row is now at row 54 (it used to be row 55). Drop it
# too; "emptyReason" and "section" shift up to rows 54 and 55.
$ws.Rows(54).Delete()

# The merged row no longer represents a named slice, so clear the leftover
# "claim" slice name from column B (Slice Name).
$ws.Range("B53").ClearContents()

# Column K ("Short") lost its longest entries (the deleted "signature" row's
# text) so its best-fit width shrinks; reflect that.
$ws.Columns("K").ColumnWidth = 84.65104166666667
